# Refresh the crypto price/volume snapshot (Price = col D, Volume(1h) = col E).
# Values are stored as plain text in the sheet (same as the scraped source),
# so each assignment is apostrophe-prefixed to force text entry and avoid
# Excel auto-converting the numeric-looking/percent strings into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.18"
$ws.Range("E2").Value = "'1.78%"
$ws.Range("D3").Value = "'27.29"
$ws.Range("E3").Value = "'1.37%"
$ws.Range("D4").Value = "'4.709"
$ws.Range("E4").Value = "'1.25%"
$ws.Range("D5").Value = "'0.06081"
$ws.Range("E5").Value = "'3.21%"
$ws.Range("D6").Value = "'6.674"
$ws.Range("E6").Value = "'0.97%"
$ws.Range("D7").Value = "'0.8460"
$ws.Range("E7").Value = "'-0.70%"
$ws.Range("D8").Value = "'0.9244"
$ws.Range("E8").Value = "'0.41%"
$ws.Range("E9").Value = "'1.97%"
$ws.Range("D10").Value = "'0.04759"
$ws.Range("E10").Value = "'13.28%"
$ws.Range("D11").Value = "'0.07101"
$ws.Range("E11").Value = "'1.51%"
$ws.Range("E12").Value = "'1.31%"
$ws.Range("D13").Value = "'0.09067"
$ws.Range("E13").Value = "'-0.45%"
$ws.Range("D14").Value = "'0.001531"
$ws.Range("E14").Value = "'-0.69%"
$ws.Range("D15").Value = "'0.0006100"
$ws.Range("E15").Value = "'0.66%"
$ws.Range("D16").Value = "'0.006182"
$ws.Range("E16").Value = "'2.35%"
$ws.Range("E17").Value = "'-0.62%"
$ws.Range("E18").Value = "'-0.87%"
$ws.Range("D19").Value = "'2.164"
$ws.Range("E19").Value = "'-0.61%"
$ws.Range("E20").Value = "'2.11%"
$ws.Range("E21").Value = "'-0.50%"
$ws.Range("D22").Value = "'4.109"
$ws.Range("E22").Value = "'4.91%"
$ws.Range("D23").Value = "'0.04247"
$ws.Range("E23").Value = "'0.23%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'-0.04%"
$ws.Range("E25").Value = "'-8.92%"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("E40").Value = "'2.32%"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("E41").Value = "'1.07%"
$ws.Range("D42").Value = "'0.004109"
$ws.Range("E42").Value = "'-34.18%"
$ws.Range("D43").Value = "'0.01636"
$ws.Range("E43").Value = "'15.63%"
$ws.Range("E44").Value = "'-8.80%"
$ws.Range("D45").Value = "'0.00005156"
$ws.Range("E45").Value = "'-3.98%"
$ws.Range("E47").Value = "'19.69%"
$ws.Range("D48").Value = "'0.1353"
$ws.Range("E48").Value = "'-46.43%"
